$d = $word.ActiveDocument

$replacements = @(
    @("2025-01-17 Friday", "2025-01-18 Saturday"),
    @("958÷5=", "321÷8="),
    @("128÷9=", "494÷4="),
    @("458÷8=", "357÷6="),
    @("943÷9=", "380÷8="),
    @("564÷4=", "788÷3="),
    @("717÷6=", "737÷4="),
    @("436÷3=", "155÷2="),
    @("525÷8=", "425÷5="),
    @("326÷5=", "586÷2="),
    @("239÷5=", "139÷2="),
    @("787÷9=", "186÷4="),
    @("199÷7=", "565÷7="),
    @("401÷2=", "979÷3="),
    @("538÷9=", "949÷4="),
    @("548÷9=", "607÷4="),
    @("355÷5=", "491÷5="),
    @("606÷9=", "680÷3="),
    @("706÷7=", "557÷8="),
    @("752÷6=", "409÷7="),
    @("850÷5=", "396÷4="),
    @("607÷8=", "729÷7="),
    @("458÷2=", "620÷8="),
    @("961÷2=", "501÷3="),
    @("662÷9=", "142÷7="),
    @("244÷4=", "624÷6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
